$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) through AC (29) hold the row payload; column A is a fixed serial index left untouched.
$firstCol = 2
$lastCol = 29

# Snapshot source rows BEFORE any writes, since several rows feed each other within a group.
$row366 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row366 += ,$ws.Cells.Item(366, $c).Value2 }
$row367 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row367 += ,$ws.Cells.Item(367, $c).Value2 }
$row368 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row368 += ,$ws.Cells.Item(368, $c).Value2 }
$row416 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row416 += ,$ws.Cells.Item(416, $c).Value2 }
$row417 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row417 += ,$ws.Cells.Item(417, $c).Value2 }
$row490 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row490 += ,$ws.Cells.Item(490, $c).Value2 }
$row491 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row491 += ,$ws.Cells.Item(491, $c).Value2 }
$row494 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row494 += ,$ws.Cells.Item(494, $c).Value2 }
$row495 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row495 += ,$ws.Cells.Item(495, $c).Value2 }
$row505 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row505 += ,$ws.Cells.Item(505, $c).Value2 }
$row506 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row506 += ,$ws.Cells.Item(506, $c).Value2 }
$row507 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row507 += ,$ws.Cells.Item(507, $c).Value2 }
$row512 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row512 += ,$ws.Cells.Item(512, $c).Value2 }
$row513 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row513 += ,$ws.Cells.Item(513, $c).Value2 }
$row519 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row519 += ,$ws.Cells.Item(519, $c).Value2 }
$row520 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row520 += ,$ws.Cells.Item(520, $c).Value2 }
$row521 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row521 += ,$ws.Cells.Item(521, $c).Value2 }
$row562 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row562 += ,$ws.Cells.Item(562, $c).Value2 }
$row563 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row563 += ,$ws.Cells.Item(563, $c).Value2 }
$row571 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row571 += ,$ws.Cells.Item(571, $c).Value2 }
$row572 = @()
for ($c = $firstCol; $c -le $lastCol; $c++) { $row572 += ,$ws.Cells.Item(572, $c).Value2 }

# Write the snapshotted rows into their new positions (rotation within each match-date group).
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(366, $c).Value = $row368[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(367, $c).Value = $row366[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(368, $c).Value = $row367[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(416, $c).Value = $row417[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(417, $c).Value = $row416[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(490, $c).Value = $row491[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(491, $c).Value = $row490[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(494, $c).Value = $row495[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(495, $c).Value = $row494[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(505, $c).Value = $row507[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(506, $c).Value = $row505[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(507, $c).Value = $row506[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(512, $c).Value = $row513[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(513, $c).Value = $row512[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(519, $c).Value = $row520[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(520, $c).Value = $row521[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(521, $c).Value = $row519[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(562, $c).Value = $row563[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(563, $c).Value = $row562[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(571, $c).Value = $row572[$c - $firstCol] }
for ($c = $firstCol; $c -le $lastCol; $c++) { $ws.Cells.Item(572, $c).Value = $row571[$c - $firstCol] }
